# "Decision_tree with embedded features"
#
# Insert a new results row for the "Decision Tree Regression+feature
# selection" model just above the existing "GB Regression" row, and update
# the trailing "Avg" row so its AVERAGE() formula keeps covering every
# model row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Push the "GB Regression" row (and the "Avg" row after it) down one slot,
# opening up row 12 for the new model.
$ws.Rows.Item(12).Insert()

# Seed the new row's formatting from the row above it (plain interior rows
# use style: bordered id / bordered name / bordered value).
$ws.Range("A11:C11").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)

# New "Decision Tree Regression+feature selection" result row.
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Decision Tree Regression+feature selection"
$ws.Range("C12").Value = 96.116233891006303

# The old row 10 ("GB Regression") is now row 11; renumber its Id.
$ws.Range("A13").Value = 11

# The "Avg" row's label cell picked up the inserted row's thick-bordered
# style via the shift; restore its original (non-filled) bordered look.
$ws.Range("B4").Copy()
$ws.Range("B14").PasteSpecial(-4122)

# Recompute the average over the now-larger C3:C13 range.
$ws.Range("C14").Formula = "=AVERAGE(C3:C13)"

# Restore the (now-shifted) active cell/selection.
$ws.Range("E9").Select() | Out-Null

$wb.Application.Calculate()
